# Insert a new weekly price-report row for "Pomelo" (Feria Lagunitas de Puerto
# Montt) above the existing row 589, pushing the rest of the historical rows
# down by one (old row 589 -> new row 590, ..., old row 710 -> new row 711).
#
# The new row reuses the same descriptive/static fields as the row that used
# to sit at 589 (market, region, product, variety, quality, unit, origin,
# kg/unit) and only carries fresh figures for: Fecha, Volumen, Precio minimo,
# Precio maximo, Precio promedio ponderado and Precio $/Kg.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 589:710 down to 590:711, leaving a blank row 589 behind.
$ws.Rows(589).Insert()

# Seed the new row with the same layout/values as the row now sitting at 590
# (the row that used to be 589 before the insert), then patch in the new
# weekly figures.
$ws.Range("A590:T590").Copy()
$ws.Range("A589").PasteSpecial()

$ws.Cells.Item(589, 4).Value = 45258    # Fecha
$ws.Cells.Item(589, 13).Value = 100     # Volumen
$ws.Cells.Item(589, 14).Value = 14000   # Precio minimo
$ws.Cells.Item(589, 15).Value = 14000   # Precio maximo
$ws.Cells.Item(589, 16).Value = 14000   # Precio promedio ponderado
$ws.Cells.Item(589, 19).Value = 1000    # Precio $/Kg
